$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Actual Result" (I) and "Status" (J) columns for the rows
# that were previously missing them. Two new messages are introduced:
#   - "Error message `n'don't forget username'" (used for rows 14 & 15)
#   - "Error message `n'username does not belong to an account'" (used for the
#     "invalid username" rows)
# The "invalid password" rows reuse the existing message from row 9
# ("Error message `n'incorrect password'").

$msgForgotUsername = "Error message `n'don't forget username'"
$msgNoAccount      = "Error message `n'username does not belong to an account'"
$msgIncorrectPwd   = "Error message `n'incorrect password'"
$pass              = "Pass"

# Row 14 (TC_LOGIN_007 - leave username & password blank)
$ws.Range("I14").Value = $msgForgotUsername

# Row 15 (TC_LOGIN_008 - blank username, invalid password)
$ws.Range("I15").Value = $msgForgotUsername
$ws.Range("J15").Value = $pass

# Row 16 (TC_LOGIN_009 - invalid username, blank password)
$ws.Range("I16").Value = $msgNoAccount
$ws.Range("J16").Value = $pass

# Row 18 (TC_LOGIN_011 - 50 char username, invalid password)
$ws.Range("I18").Value = $msgIncorrectPwd
$ws.Range("J18").Value = $pass

# Row 20 (TC_LOGIN_013 - invalid username, 50 char password)
$ws.Range("I20").Value = $msgNoAccount
$ws.Range("J20").Value = $pass

# Row 22 (TC_LOGIN_015 - 100 char username, invalid password)
$ws.Range("I22").Value = $msgIncorrectPwd
$ws.Range("J22").Value = $pass

# Row 24 (TC_LOGIN_017 - invalid username, 100 char password)
$ws.Range("I24").Value = $msgNoAccount
$ws.Range("J24").Value = $pass

# Row 26 (TC_LOGIN_019 - 128 char username, invalid password)
$ws.Range("I26").Value = $msgIncorrectPwd
$ws.Range("J26").Value = $pass

# Row 28 (TC_LOGIN_021 - invalid username, 128 char password)
$ws.Range("I28").Value = $msgNoAccount
$ws.Range("J28").Value = $pass

# Row 30 (TC_LOGIN_023 - 200 char username, invalid password)
$ws.Range("I30").Value = $msgIncorrectPwd
$ws.Range("J30").Value = $pass

# Row 32 (TC_LOGIN_025 - invalid username, 200 char password)
$ws.Range("I32").Value = $msgNoAccount
$ws.Range("J32").Value = $pass

# --- Remove the two obsolete test cases (TC_LOGIN_026 and TC_LOGIN_027),
# which occupy rows 33 and 34. Deleting these rows shifts the trailing
# "  " row (35) up to row 33, and shrinks the Table3 / dimension / data
# validation ranges automatically.
$ws.Range("A33:A34").EntireRow.Delete() | Out-Null

# --- Update the view: scroll so column D is left-most visible and select
# the last two rows' Status cells, landing on J32.
$ws.Range("J30,J32").Select() | Out-Null
$ws.Range("J32").Activate() | Out-Null
